$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts B..L to C..M)
$ws.Columns("B:B").Insert()

# New column B holds a running index 0..49 for rows 2..51
$ws.Range("B2").Value = 0
$ws.Range("B3").Formula = "=B2+1"
$ws.Range("B4:B51").FormulaR1C1 = "=R[-1]C+1"

# Bold style for entire column B (like the header column style)
$ws.Range("B2:B51").Font.Bold = $true

# Column B width
$ws.Columns("B:B").ColumnWidth = 7.44140625

# Update selection to match the target view
$ws.Range("C7").Select()
